$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Target" cash-in (I20) is now a fixed, hard-coded figure instead of
# being derived from the margin (I21); replaces the old formula.
$ws.Range("I20").Value = 150000

# "Margin" (I21) is now computed FROM the target/profits relationship
# instead of being a manually typed assumption, and is shown with one
# decimal place.
$ws.Range("I21").Formula = "=(I20-O17)/O17"
$ws.Range("I21").NumberFormat = "0.0%"

# Update the active selection to the cell the author was last working on.
$ws.Range("I21").Select()
